# Scheduled runner update: refresh Universalis market price & profit data
# across the Leve profit-tracking sheets (columns H-N).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 298.6154
$ws.Range("I96").Value = 286.66666
$ws.Range("J96").Value = 325.5
$ws.Range("K96").Value = 859.9999799999999
$ws.Range("L96").Value = 976.5
$ws.Range("M96").Value = 513.0000200000001
$ws.Range("N96").Value = -3722.5

$ws.Range("H137").Value = 35715896
$ws.Range("I137").Value = 1323.8948
$ws.Range("J137").Value = 111113336
$ws.Range("K137").Value = 3971.6844
$ws.Range("L137").Value = 333340008
$ws.Range("M137").Value = -1421.6844
$ws.Range("N137").Value = -333345108

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8165.7017
$ws.Range("I32").Value = 6924.2256
$ws.Range("J32").Value = 23560
$ws.Range("K32").Value = 6924.2256
$ws.Range("L32").Value = 23560
$ws.Range("M32").Value = -6637.2256
$ws.Range("N32").Value = -24134

$ws.Range("H37").Value = 9522.111000000001
$ws.Range("J37").Value = 10087.375
$ws.Range("L37").Value = 10087.375
$ws.Range("N37").Value = -10633.375

$ws.Range("H74").Value = 4687.9165
$ws.Range("I74").Value = 739.625
$ws.Range("J74").Value = 12584.5
$ws.Range("K74").Value = 739.625
$ws.Range("L74").Value = 12584.5
$ws.Range("M74").Value = 134.375
$ws.Range("N74").Value = -14332.5

$ws.Range("H77").Value = 4687.9165
$ws.Range("I77").Value = 739.625
$ws.Range("J77").Value = 12584.5
$ws.Range("K77").Value = 3698.125
$ws.Range("L77").Value = 62922.5
$ws.Range("M77").Value = 669.875
$ws.Range("N77").Value = -71658.5

$ws.Range("H132").Value = 45467.25
$ws.Range("I132").Value = 3960.8
$ws.Range("J132").Value = 252999.5
$ws.Range("K132").Value = 11882.4
$ws.Range("L132").Value = 758998.5
$ws.Range("M132").Value = -9352.400000000001
$ws.Range("N132").Value = -764058.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 66066.7
$ws.Range("I134").Value = 85910
$ws.Range("K134").Value = 257730
$ws.Range("M134").Value = -255195

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1364.0333
$ws.Range("I31").Value = 788.7917
$ws.Range("J31").Value = 3665
$ws.Range("K31").Value = 788.7917
$ws.Range("L31").Value = 3665
$ws.Range("M31").Value = -493.7917
$ws.Range("N31").Value = -4255

$ws.Range("H34").Value = 1364.0333
$ws.Range("I34").Value = 788.7917
$ws.Range("J34").Value = 3665
$ws.Range("K34").Value = 788.7917
$ws.Range("L34").Value = 3665
$ws.Range("M34").Value = -586.7917
$ws.Range("N34").Value = -4069

$ws.Range("H51").Value = 9667
$ws.Range("J51").Value = 10271.857
$ws.Range("L51").Value = 10271.857
$ws.Range("N51").Value = -11743.857

$ws.Range("H59").Value = 16618.455
$ws.Range("J59").Value = 16680.3
$ws.Range("L59").Value = 16680.3
$ws.Range("N59").Value = -18970.3

$ws.Range("H60").Value = 9518.727999999999
$ws.Range("J60").Value = 9970.6
$ws.Range("L60").Value = 9970.6
$ws.Range("N60").Value = -10992.6

$ws.Range("H61").Value = 9667
$ws.Range("J61").Value = 10271.857
$ws.Range("L61").Value = 10271.857
$ws.Range("N61").Value = -10967.857

$ws.Range("H68").Value = 18889.223
$ws.Range("J68").Value = 18889.223
$ws.Range("L68").Value = 18889.223
$ws.Range("N68").Value = -20387.223

$ws.Range("H71").Value = 18889.223
$ws.Range("J71").Value = 18889.223
$ws.Range("L71").Value = 56667.66900000001
$ws.Range("N71").Value = -64155.66900000001

$ws.Range("H74").Value = 14967.6
$ws.Range("J74").Value = 16186.223
$ws.Range("L74").Value = 16186.223
$ws.Range("N74").Value = -17934.223

$ws.Range("H77").Value = 14967.6
$ws.Range("J77").Value = 16186.223
$ws.Range("L77").Value = 48558.669
$ws.Range("N77").Value = -57294.669

$ws.Range("H134").Value = 1929.317
$ws.Range("I134").Value = 1815.6129
$ws.Range("J134").Value = 2281.8
$ws.Range("K134").Value = 5446.8387
$ws.Range("L134").Value = 6845.400000000001
$ws.Range("M134").Value = -2911.8387
$ws.Range("N134").Value = -11915.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1396
$ws.Range("J22").Value = 1497.5
$ws.Range("L22").Value = 4492.5
$ws.Range("N22").Value = -4830.5

$ws.Range("H27").Value = 1396
$ws.Range("J27").Value = 1497.5
$ws.Range("L27").Value = 4492.5
$ws.Range("N27").Value = -4696.5

$ws.Range("H131").Value = 4069.6462
$ws.Range("I131").Value = 5439.9165
$ws.Range("J131").Value = 3759.3962
$ws.Range("K131").Value = 16319.7495
$ws.Range("L131").Value = 11278.1886
$ws.Range("M131").Value = -11279.7495
$ws.Range("N131").Value = -21358.1886

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 23713.8
$ws.Range("J32").Value = 23713.8
$ws.Range("L32").Value = 23713.8
$ws.Range("N32").Value = -24305.8

$ws.Range("H111").Value = 24664.334
$ws.Range("J111").Value = 24664.334
$ws.Range("L111").Value = 24664.334
$ws.Range("N111").Value = -30798.334

$ws.Range("H122").Value = 3808.2104
$ws.Range("I122").Value = 4144.4707
$ws.Range("J122").Value = 950
$ws.Range("K122").Value = 12433.4121
$ws.Range("L122").Value = 2850
$ws.Range("M122").Value = -9983.4121
$ws.Range("N122").Value = -7750

$ws.Range("H132").Value = 3285
$ws.Range("I132").Value = 2837.4285
$ws.Range("J132").Value = 4329.3335
$ws.Range("K132").Value = 8512.2855
$ws.Range("L132").Value = 12988.0005
$ws.Range("M132").Value = -5982.2855
$ws.Range("N132").Value = -18048.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 15340
$ws.Range("J74").Value = 16800
$ws.Range("L74").Value = 16800
$ws.Range("N74").Value = -18672

$ws.Range("H77").Value = 15340
$ws.Range("J77").Value = 16800
$ws.Range("L77").Value = 50400
$ws.Range("N77").Value = -59760

$ws.Range("H107").Value = 479.33334
$ws.Range("I107").Value = 284.81818
$ws.Range("K107").Value = 854.45454
$ws.Range("M107").Value = 1065.54546

$ws.Range("H136").Value = 6747.381
$ws.Range("I136").Value = 9132.861999999999
$ws.Range("J136").Value = 1425.9231
$ws.Range("K136").Value = 27398.586
$ws.Range("L136").Value = 4277.7693
$ws.Range("M136").Value = -24848.586
$ws.Range("N136").Value = -9377.7693
